$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (B1, C1, D1)
$ws.Range("B1").Value = "Remained Occupied"
$ws.Range("C1").Value = "Extirpated"
$ws.Range("D1").Value = "Colonized"

# Update data rows: B becomes numeric, C and D get new values
$data = @{
    2 = @(1130, 33, 41)
    3 = @(1874, 142, 80)
    4 = @(1642, 126, 112)
    5 = @(684, 90, 39)
    6 = @(561, 18, 21)
    7 = @(383, 18, 40)
    8 = @(77, 8, 4)
    9 = @(110, 2, 6)
    10 = @(62, 2, 3)
    11 = @(965, 55, 48)
    12 = @(985, 74, 71)
    13 = @(1240, 89, 62)
    14 = @(658, 91, 46)
    15 = @(233, 36, 28)
    16 = @(173, 30, 8)
    17 = @(61, 22, 16)
    18 = @(2058, 198, 135)
    19 = @(608, 29, 14)
    20 = @(876, 115, 38)
    21 = @(538, 66, 67)
    22 = @(4374, 207, 232)
    23 = @(171, 24, 14)
    24 = @(190, 13, 3)
    25 = @(7898, 371, 360)
    26 = @(238, 2, 4)
    27 = @(195, 25, 11)
    28 = @(1388, 73, 99)
    29 = @(2585, 160, 203)
    30 = @(724, 19, 30)
    31 = @(853, 88, 88)
    32 = @(69, 10, 5)
    33 = @(1128, 190, 158)
    34 = @(63, 8, 3)
    35 = @(754, 85, 49)
    36 = @(66, 14, 4)
    37 = @(289, 56, 39)
    38 = @(66, 29, 18)
    39 = @(54, 3, 11)
    40 = @(611, 39, 32)
    41 = @(98, 21, 8)
    42 = @(149, 19, 18)
    43 = @(188, 11, 6)
    44 = @(791, 60, 30)
    45 = @(335, 5, 0)
    46 = @(321, 13, 11)
    47 = @(805, 91, 42)
    48 = @(60, 2, 0)
    49 = @(184, 23, 13)
    50 = @(256, 8, 21)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
